$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing F/G values for rows with revised Ag test figures
$ws.Range("F630").Value = 46872

$ws.Range("F632").Value = 44388
$ws.Range("G632").Value = 2669

$ws.Range("F635").Value = 83636

$ws.Range("F637").Value = 43767
$ws.Range("G637").Value = 2117

$ws.Range("F639").Value = 40732

$ws.Range("F642").Value = 67487
$ws.Range("G642").Value = 2384

$ws.Range("F643").Value = 43475

$ws.Range("F645").Value = 35578
$ws.Range("G645").Value = 1307

$ws.Range("F649").Value = 62465
$ws.Range("G649").Value = 1809

$ws.Range("F650").Value = 38063

$ws.Range("F651").Value = 37127

$ws.Range("F652").Value = 35048
$ws.Range("G652").Value = 1094

$ws.Range("F653").Value = 34144

$ws.Range("F656").Value = 52372
$ws.Range("G656").Value = 1240

$ws.Range("F658").Value = 27157

$ws.Range("F659").Value = 26199

$ws.Range("F662").Value = 12743
$ws.Range("G662").Value = 541

$ws.Range("F663").Value = 37071

$ws.Range("F664").Value = 26402
$ws.Range("G664").Value = 773

$ws.Range("F665").Value = 24300

$ws.Range("F666").Value = 23747

$ws.Range("F667").Value = 17163
$ws.Range("G667").Value = 608

$ws.Range("F668").Value = 3369

$ws.Range("F669").Value = 23378
$ws.Range("G669").Value = 634

$ws.Range("F670").Value = 52651
$ws.Range("G670").Value = 911

$ws.Range("F671").Value = 33018
$ws.Range("G671").Value = 617

$ws.Range("F672").Value = 29555
$ws.Range("G672").Value = 578

$ws.Range("F673").Value = 10054
$ws.Range("G673").Value = 307

$ws.Range("F674").Value = 28185
$ws.Range("G674").Value = 677

$ws.Range("F675").Value = 13315
$ws.Range("G675").Value = 352

$ws.Range("F676").Value = 27604
$ws.Range("G676").Value = 436

$ws.Range("F677").Value = 53762
$ws.Range("G677").Value = 781

# Row 678 was missing F/G values; add them now
$ws.Range("F678").Value = 31998
$ws.Range("G678").Value = 611

# Append new row 679 with the latest day's data
$ws.Range("A679").Value = 44573
$ws.Range("A679").NumberFormat = $ws.Range("A678").NumberFormat
$ws.Range("B679").Value = 872511
$ws.Range("C679").Value = 11680
$ws.Range("D679").Value = 2848
$ws.Range("E679").Value = 17128
$ws.Range("F679").Value = 18271
$ws.Range("G679").Value = 389
